# Adição da descrição dos entregaveis no backlog
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HFSystem")

$ws.Range("C4").Value  = "LOCAÇÃO DOS ELEMENTOS DO PROJETO EM UM REPOSITÓRIO NO GITHUB"
$ws.Range("C12").Value = "CONFIGURAÇÃO DE TODOS OS ENTREGÁVEIS NO TRELLO"
$ws.Range("C13").Value = "CRIAÇÃO DO BACKLOG"
$ws.Range("C20").Value = "CRIAÇÃO DE PLANILHA DE RISCOS JUNTAMENTE COM GRÁFICO DE BURNDOWN"
$ws.Range("C21").Value = "CRIAÇÃO DE UMA DASHBOARD ESTÁTICA"
$ws.Range("C28").Value = "MODELAGEM DO BANCO DE DADOS"
$ws.Range("C29").Value = "UTILIZAÇÃO DO SENSOR JUNTO COM A API"
$ws.Range("C34").Value = "INTEGRAÇÃO DO ARDUINO COM O BANCO DE DADOS"
$ws.Range("C35").Value = "A DEFINIR "
$ws.Range("C36").Value = "SINCRONIZAÇÃO DOS DADOS CAPTADOS PELO SENSOR COM A API"
$ws.Range("C37").Value = "A DEFINIR "
$ws.Range("C38").Value = "A DEFINIR"
$ws.Range("C39").Value = "CRIAÇÃO DE MODELAGEM DEFINITIVA DO PROJETO"
$ws.Range("C40").Value = "A DEFINIR"
$ws.Range("C41").Value = "A DEFINIR"
$ws.Range("C42").Value = "A DEFINIR"
$ws.Range("C43").Value = "CRIAÇÃO DE SITE INSTITUCIONAL EM HTML/CSS/JS FUNCIONAL"
$ws.Range("C44").Value = "CRIAÇÃO E CONFIGURAÇÃO DE UMA DASHBOARD FUNCIONAL"
$ws.Range("C45").Value = "CRIAÇÃO DE PÁGINA DE LOGIN E CADASTRO FUNCIONAL"
$ws.Range("C46").Value = "A DEFINIR"
